# Refresh market-price-driven profit columns (H:N) across all 8 job sheets,
# per the scheduled market-data runner. The workbook caches plain numbers
# (no live formulas in H:N), so each target cell is written directly with
# its freshly-pulled value. A few cells had no result this run and are
# cleared instead of being set to a number.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 633.8  # H2: was 809.8570999999999
$ws.Cells.Item(2, 9).Value = 317.25  # I2: was 353.8
$ws.Cells.Item(2, 10).Value = 1900  # J2: was 1950
$ws.Cells.Item(2, 11).Value = 317.25  # K2: was 353.8
$ws.Cells.Item(2, 12).Value = 1900  # L2: was 1950
$ws.Cells.Item(2, 13).Value = -204.25  # M2: was -240.8
$ws.Cells.Item(2, 14).Value = -2126  # N2: was -2176
$ws.Cells.Item(17, 8).Value = 7279.263  # H17: was 7968.0586
$ws.Cells.Item(17, 10).Value = 7583.6665  # J17: was 8353.5625
$ws.Cells.Item(17, 12).Value = 22750.9995  # L17: was 25060.6875
$ws.Cells.Item(17, 14).Value = -23086.9995  # N17: was -25396.6875
$ws.Cells.Item(19, 8).Value = 1091.6923  # H19: was 1017.3571
$ws.Cells.Item(19, 9).Value = 219.75  # I19: was 186
$ws.Cells.Item(19, 11).Value = 219.75  # K19: was 186
$ws.Cells.Item(19, 13).Value = -44.75  # M19: was -11
$ws.Cells.Item(21, 8).Value = 54281.9  # H21: was 13356.333
$ws.Cells.Item(21, 9).Value = 5000  # I21: was 5025
$ws.Cells.Item(21, 10).Value = 59757.668  # J21: was 30019
$ws.Cells.Item(21, 11).Value = 5000  # K21: was 5025
$ws.Cells.Item(21, 12).Value = 59757.668  # L21: was 30019
$ws.Cells.Item(21, 13).Value = -4532  # M21: was -4557
$ws.Cells.Item(21, 14).Value = -60693.668  # N21: was -30955
$ws.Cells.Item(23, 8).Value = 54281.9  # H23: was 13356.333
$ws.Cells.Item(23, 9).Value = 5000  # I23: was 5025
$ws.Cells.Item(23, 10).Value = 59757.668  # J23: was 30019
$ws.Cells.Item(23, 11).Value = 5000  # K23: was 5025
$ws.Cells.Item(23, 12).Value = 59757.668  # L23: was 30019
$ws.Cells.Item(23, 13).Value = -4766  # M23: was -4791
$ws.Cells.Item(23, 14).Value = -60225.668  # N23: was -30487
$ws.Cells.Item(38, 8).Value = 2749.3  # H38: was 1728.5454
$ws.Cells.Item(38, 9).Value = 71.57143000000001  # I38: was 64.875
$ws.Cells.Item(38, 10).Value = 8997.333000000001  # J38: was 6165
$ws.Cells.Item(38, 11).Value = 214.71429  # K38: was 194.625
$ws.Cells.Item(38, 12).Value = 26991.999  # L38: was 18495
$ws.Cells.Item(38, 13).Value = 157.28571  # M38: was 177.375
$ws.Cells.Item(38, 14).Value = -27735.999  # N38: was -19239
$ws.Cells.Item(58, 8).Value = 577.8  # H58: was 575
$ws.Cells.Item(58, 9).Value = 222.25  # I58: was 158.33333
$ws.Cells.Item(58, 10).Value = 2000  # J58: was 1200
$ws.Cells.Item(58, 11).Value = 666.75  # K58: was 474.99999
$ws.Cells.Item(58, 12).Value = 6000  # L58: was 3600
$ws.Cells.Item(58, 13).Value = -516.75  # M58: was -324.99999
$ws.Cells.Item(58, 14).Value = -6300  # N58: was -3900
$ws.Cells.Item(64, 8).Value = 7854.25  # H64: was 7657.483
$ws.Cells.Item(64, 9).Value = 5045.75  # I64: was 5070.75
$ws.Cells.Item(64, 10).Value = 8415.950000000001  # J64: was 8071.36
$ws.Cells.Item(64, 11).Value = 5045.75  # K64: was 5070.75
$ws.Cells.Item(64, 12).Value = 8415.950000000001  # L64: was 8071.36
$ws.Cells.Item(64, 13).Value = -4797.75  # M64: was -4822.75
$ws.Cells.Item(64, 14).Value = -8911.950000000001  # N64: was -8567.360000000001
$ws.Cells.Item(67, 8).Value = 7854.25  # H67: was 7657.483
$ws.Cells.Item(67, 9).Value = 5045.75  # I67: was 5070.75
$ws.Cells.Item(67, 10).Value = 8415.950000000001  # J67: was 8071.36
$ws.Cells.Item(67, 11).Value = 5045.75  # K67: was 5070.75
$ws.Cells.Item(67, 12).Value = 8415.950000000001  # L67: was 8071.36
$ws.Cells.Item(67, 13).Value = -4187.75  # M67: was -4212.75
$ws.Cells.Item(67, 14).Value = -10131.95  # N67: was -9787.360000000001
$ws.Cells.Item(74, 8).Value = 6248.4595  # H74: was 6647.6665
$ws.Cells.Item(74, 9).Value = 4298.4546  # I74: was 5326.143
$ws.Cells.Item(74, 10).Value = 7073.4614  # J74: was 7488.636
$ws.Cells.Item(74, 11).Value = 4298.4546  # K74: was 5326.143
$ws.Cells.Item(74, 12).Value = 7073.4614  # L74: was 7488.636
$ws.Cells.Item(74, 13).Value = -3362.4546  # M74: was -4390.143
$ws.Cells.Item(74, 14).Value = -8945.4614  # N74: was -9360.636
$ws.Cells.Item(77, 8).Value = 6248.4595  # H77: was 6647.6665
$ws.Cells.Item(77, 9).Value = 4298.4546  # I77: was 5326.143
$ws.Cells.Item(77, 10).Value = 7073.4614  # J77: was 7488.636
$ws.Cells.Item(77, 11).Value = 21492.273  # K77: was 26630.715
$ws.Cells.Item(77, 12).Value = 35367.307  # L77: was 37443.18
$ws.Cells.Item(77, 13).Value = -16812.273  # M77: was -21950.715
$ws.Cells.Item(77, 14).Value = -44727.307  # N77: was -46803.18
$ws.Cells.Item(116, 8).Value = 5529.6665  # H116: was 5392.7407
$ws.Cells.Item(116, 9).Value = 3383.375  # I116: was 3251.6667
$ws.Cells.Item(116, 10).Value = 6433.3687  # J116: was 6463.278
$ws.Cells.Item(116, 11).Value = 3383.375  # K116: was 3251.6667
$ws.Cells.Item(116, 12).Value = 6433.3687  # L116: was 6463.278
$ws.Cells.Item(116, 13).Value = 58.625  # M116: was 190.3332999999998
$ws.Cells.Item(116, 14).Value = -13317.3687  # N116: was -13347.278
$ws.Cells.Item(120, 8).Value = 80000  # H120: was 71424.28999999999
$ws.Cells.Item(120, 10).Value = 80000  # J120: was 71424.28999999999
$ws.Cells.Item(120, 12).Value = 80000  # L120: was 71424.28999999999
$ws.Cells.Item(120, 14).Value = -89676  # N120: was -81100.28999999999
$ws.Cells.Item(132, 8).Value = 59102.832  # H132: was 60788.63
$ws.Cells.Item(132, 9).Value = 60705.914  # I132: was 62488.44
$ws.Cells.Item(132, 11).Value = 182117.742  # K132: was 187465.32
$ws.Cells.Item(132, 13).Value = -179587.742  # M132: was -184935.32

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7682.577  # H32: was 7689.4614
$ws.Cells.Item(32, 9).Value = 4740.661  # I32: was 4749.3228
$ws.Cells.Item(32, 11).Value = 4740.661  # K32: was 4749.3228
$ws.Cells.Item(32, 13).Value = -4453.661  # M32: was -4462.3228
$ws.Cells.Item(74, 8).Value = 150972.84  # H74: was 110073.89
$ws.Cells.Item(74, 9).Value = 129968.625  # I74: was 87701.5
$ws.Cells.Item(74, 10).Value = 184579.6  # J74: was 154818.67
$ws.Cells.Item(74, 11).Value = 129968.625  # K74: was 87701.5
$ws.Cells.Item(74, 12).Value = 184579.6  # L74: was 154818.67
$ws.Cells.Item(74, 13).Value = -129094.625  # M74: was -86827.5
$ws.Cells.Item(74, 14).Value = -186327.6  # N74: was -156566.67
$ws.Cells.Item(77, 8).Value = 150972.84  # H77: was 110073.89
$ws.Cells.Item(77, 9).Value = 129968.625  # I77: was 87701.5
$ws.Cells.Item(77, 10).Value = 184579.6  # J77: was 154818.67
$ws.Cells.Item(77, 11).Value = 649843.125  # K77: was 438507.5
$ws.Cells.Item(77, 12).Value = 922898  # L77: was 774093.3500000001
$ws.Cells.Item(77, 13).Value = -645475.125  # M77: was -434139.5
$ws.Cells.Item(77, 14).Value = -931634  # N77: was -782829.3500000001
$ws.Cells.Item(98, 8).Value = 0  # H98: was 25000
$ws.Cells.Item(98, 10).Value = 0  # J98: was 25000
$ws.Cells.Item(98, 12).Value = 0  # L98: was 25000
$ws.Cells.Item(98, 14).ClearContents()  # N98: was -30990, now blank
$ws.Cells.Item(122, 8).Value = 3343563.8  # H122: was 3677877.2
$ws.Cells.Item(122, 9).Value = 3762149.8  # I122: was 4388854
$ws.Cells.Item(122, 10).Value = 2611038  # J122: was 2611412.5
$ws.Cells.Item(122, 11).Value = 11286449.4  # K122: was 13166562
$ws.Cells.Item(122, 12).Value = 7833114  # L122: was 7834237.5
$ws.Cells.Item(122, 13).Value = -11283999.4  # M122: was -13164112
$ws.Cells.Item(122, 14).Value = -7838014  # N122: was -7839137.5
$ws.Cells.Item(132, 8).Value = 4580.609  # H132: was 4478.617
$ws.Cells.Item(132, 9).Value = 5259  # I132: was 4946.8213
$ws.Cells.Item(132, 10).Value = 3698.7  # J132: was 3788.6316
$ws.Cells.Item(132, 11).Value = 15777  # K132: was 14840.4639
$ws.Cells.Item(132, 12).Value = 11096.1  # L132: was 11365.8948
$ws.Cells.Item(132, 13).Value = -13247  # M132: was -12310.4639
$ws.Cells.Item(132, 14).Value = -16156.1  # N132: was -16425.8948

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2385672.5  # H86: was 2443842.5
$ws.Cells.Item(86, 9).Value = 3127096  # I86: was 3227947.5
$ws.Cells.Item(86, 11).Value = 3127096  # K86: was 3227947.5
$ws.Cells.Item(86, 13).Value = -3125973  # M86: was -3226824.5
$ws.Cells.Item(89, 8).Value = 2385672.5  # H89: was 2443842.5
$ws.Cells.Item(89, 9).Value = 3127096  # I89: was 3227947.5
$ws.Cells.Item(89, 11).Value = 15635480  # K89: was 16139737.5
$ws.Cells.Item(89, 13).Value = -15629864  # M89: was -16134121.5
$ws.Cells.Item(100, 8).Value = 21185.8  # H100: was 30000
$ws.Cells.Item(100, 10).Value = 21185.8  # J100: was 30000
$ws.Cells.Item(100, 12).Value = 21185.8  # L100: was 30000
$ws.Cells.Item(100, 14).Value = -23349.8  # N100: was -32164
$ws.Cells.Item(105, 8).Value = 3473395.2  # H105: was 3677677.2
$ws.Cells.Item(105, 9).Value = 3677536  # I105: was 3907344.5
$ws.Cells.Item(105, 11).Value = 3677536  # K105: was 3907344.5
$ws.Cells.Item(105, 13).Value = -3675789  # M105: was -3905597.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 11553.637  # H86: was 7668.8887
$ws.Cells.Item(86, 9).Value = 10093.742  # I86: was 6422.5557
$ws.Cells.Item(86, 10).Value = 15034.923  # J86: was 10161.556
$ws.Cells.Item(86, 11).Value = 10093.742  # K86: was 6422.5557
$ws.Cells.Item(86, 12).Value = 15034.923  # L86: was 10161.556
$ws.Cells.Item(86, 13).Value = -8970.742  # M86: was -5299.5557
$ws.Cells.Item(86, 14).Value = -17280.923  # N86: was -12407.556
$ws.Cells.Item(89, 8).Value = 11553.637  # H89: was 7668.8887
$ws.Cells.Item(89, 9).Value = 10093.742  # I89: was 6422.5557
$ws.Cells.Item(89, 10).Value = 15034.923  # J89: was 10161.556
$ws.Cells.Item(89, 11).Value = 50468.71  # K89: was 32112.7785
$ws.Cells.Item(89, 12).Value = 75174.61500000001  # L89: was 50807.78
$ws.Cells.Item(89, 13).Value = -44852.71  # M89: was -26496.7785
$ws.Cells.Item(89, 14).Value = -86406.61500000001  # N89: was -62039.78
$ws.Cells.Item(132, 8).Value = 120014.69  # H132: was 127935.266
$ws.Cells.Item(132, 9).Value = 92473.37  # I132: was 101600.1
$ws.Cells.Item(132, 11).Value = 277420.11  # K132: was 304800.3
$ws.Cells.Item(132, 13).Value = -274890.11  # M132: was -302270.3
$ws.Cells.Item(141, 8).Value = 196095.77  # H141: was 228205
$ws.Cells.Item(141, 10).Value = 211579.08  # J141: was 249995.9
$ws.Cells.Item(141, 12).Value = 211579.08  # L141: was 249995.9
$ws.Cells.Item(141, 14).Value = -221939.08  # N141: was -260355.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 77773.92  # H5: was 67459.664
$ws.Cells.Item(5, 9).Value = 778.36365  # I5: was 722.7692
$ws.Cells.Item(5, 11).Value = 2335.09095  # K5: was 2168.3076
$ws.Cells.Item(5, 13).Value = -2223.09095  # M5: was -2056.3076
$ws.Cells.Item(59, 8).Value = 1000  # H59: was 2000
$ws.Cells.Item(59, 9).Value = 1000  # I59: was 2000
$ws.Cells.Item(59, 11).Value = 3000  # K59: was 6000
$ws.Cells.Item(59, 13).Value = -2460  # M59: was -5460
$ws.Cells.Item(135, 8).Value = 77773.92  # H135: was 67459.664
$ws.Cells.Item(135, 9).Value = 778.36365  # I135: was 722.7692
$ws.Cells.Item(135, 11).Value = 7005.27285  # K135: was 6504.922799999999
$ws.Cells.Item(135, 13).Value = -4470.27285  # M135: was -3969.922799999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 18185522  # H70: was 33336872
$ws.Cells.Item(70, 9).Value = 22225836  # I70: was 50003252
$ws.Cells.Item(70, 11).Value = 22225836  # K70: was 50003252
$ws.Cells.Item(70, 13).Value = -22225566  # M70: was -50002982
$ws.Cells.Item(73, 8).Value = 18185522  # H73: was 33336872
$ws.Cells.Item(73, 9).Value = 22225836  # I73: was 50003252
$ws.Cells.Item(73, 11).Value = 22225836  # K73: was 50003252
$ws.Cells.Item(73, 13).Value = -22224900  # M73: was -50002316
$ws.Cells.Item(122, 8).Value = 265524.9  # H122: was 290969.03
$ws.Cells.Item(122, 9).Value = 358591.44  # I122: was 389651.6
$ws.Cells.Item(122, 10).Value = 7006.778  # J122: was 7256.625
$ws.Cells.Item(122, 11).Value = 1075774.32  # K122: was 1168954.8
$ws.Cells.Item(122, 12).Value = 21020.334  # L122: was 21769.875
$ws.Cells.Item(122, 13).Value = -1073324.32  # M122: was -1166504.8
$ws.Cells.Item(122, 14).Value = -25920.334  # N122: was -26669.875
$ws.Cells.Item(132, 8).Value = 8345.527  # H132: was 8363.666999999999
$ws.Cells.Item(132, 9).Value = 6158.552  # I132: was 6161.0347
$ws.Cells.Item(132, 10).Value = 17405.857  # J132: was 17488.857
$ws.Cells.Item(132, 11).Value = 18475.656  # K132: was 18483.1041
$ws.Cells.Item(132, 12).Value = 52217.571  # L132: was 52466.571
$ws.Cells.Item(132, 13).Value = -15945.656  # M132: was -15953.1041
$ws.Cells.Item(132, 14).Value = -57277.571  # N132: was -57526.571

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 37213.48  # H22: was 44134.383
$ws.Cells.Item(22, 9).Value = 89745.7  # I22: was 111979
$ws.Cells.Item(22, 10).Value = 2192  # J22: was 2383.8462
$ws.Cells.Item(22, 11).Value = 89745.7  # K22: was 111979
$ws.Cells.Item(22, 12).Value = 2192  # L22: was 2383.8462
$ws.Cells.Item(22, 13).Value = -89450.7  # M22: was -111684
$ws.Cells.Item(22, 14).Value = -2782  # N22: was -2973.8462
$ws.Cells.Item(27, 8).Value = 37213.48  # H27: was 44134.383
$ws.Cells.Item(27, 9).Value = 89745.7  # I27: was 111979
$ws.Cells.Item(27, 10).Value = 2192  # J27: was 2383.8462
$ws.Cells.Item(27, 11).Value = 89745.7  # K27: was 111979
$ws.Cells.Item(27, 12).Value = 2192  # L27: was 2383.8462
$ws.Cells.Item(27, 13).Value = -89638.7  # M27: was -111872
$ws.Cells.Item(27, 14).Value = -2406  # N27: was -2597.8462
$ws.Cells.Item(38, 8).Value = 1000000000  # H38: was 500007000
$ws.Cells.Item(38, 10).Value = 0  # J38: was 14000
$ws.Cells.Item(38, 12).Value = 0  # L38: was 14000
$ws.Cells.Item(38, 14).ClearContents()  # N38: was -14820, now blank
$ws.Cells.Item(40, 8).Value = 12474.5  # H40: was 7193.615
$ws.Cells.Item(40, 9).Value = 0  # I40: was 4846.5557
$ws.Cells.Item(40, 11).Value = 0  # K40: was 4846.5557
$ws.Cells.Item(40, 13).ClearContents()  # M40: was -4710.5557, now blank
$ws.Cells.Item(55, 8).Value = 1053.7727  # H55: was 1107.9524
$ws.Cells.Item(55, 9).Value = 919.9091  # I55: was 1004.1
$ws.Cells.Item(55, 10).Value = 1187.6364  # J55: was 1202.3636
$ws.Cells.Item(55, 11).Value = 919.9091  # K55: was 1004.1
$ws.Cells.Item(55, 12).Value = 1187.6364  # L55: was 1202.3636
$ws.Cells.Item(55, 13).Value = -746.9091  # M55: was -831.1
$ws.Cells.Item(55, 14).Value = -1533.6364  # N55: was -1548.3636
$ws.Cells.Item(56, 8).Value = 17274.75  # H56: was 10433.333
$ws.Cells.Item(56, 10).Value = 19900  # J56: was 2001
$ws.Cells.Item(56, 12).Value = 19900  # L56: was 2001
$ws.Cells.Item(56, 14).Value = -21282  # N56: was -3383
$ws.Cells.Item(122, 8).Value = 9471.75  # H122: was 7497.0713
$ws.Cells.Item(122, 9).Value = 4995  # I122: was 4399
$ws.Cells.Item(122, 10).Value = 10111.286  # J122: was 9218.223
$ws.Cells.Item(122, 11).Value = 14985  # K122: was 13197
$ws.Cells.Item(122, 12).Value = 30333.858  # L122: was 27654.669
$ws.Cells.Item(122, 13).Value = -12535  # M122: was -10747
$ws.Cells.Item(122, 14).Value = -35233.858  # N122: was -32554.669
$ws.Cells.Item(132, 8).Value = 10624.2  # H132: was 10702.86
$ws.Cells.Item(132, 9).Value = 11867.1  # I132: was 11904.225
$ws.Cells.Item(132, 10).Value = 5652.6  # J132: was 5897.4
$ws.Cells.Item(132, 11).Value = 35601.3  # K132: was 35712.675
$ws.Cells.Item(132, 12).Value = 16957.8  # L132: was 17692.2
$ws.Cells.Item(132, 13).Value = -33071.3  # M132: was -33182.675
$ws.Cells.Item(132, 14).Value = -22017.8  # N132: was -22752.2
$ws.Cells.Item(136, 8).Value = 60644.805  # H136: was 56159.31
$ws.Cells.Item(136, 9).Value = 170311.33  # I136: was 136715.73
$ws.Cells.Item(136, 11).Value = 510933.99  # K136: was 410147.1900000001
$ws.Cells.Item(136, 13).Value = -508383.99  # M136: was -407597.1900000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 2023.25  # H113: was 2098.182
$ws.Cells.Item(113, 9).Value = 1347  # I113: was 1495
$ws.Cells.Item(113, 11).Value = 4041  # K113: was 4485
$ws.Cells.Item(113, 13).Value = -1871  # M113: was -2315
$ws.Cells.Item(122, 8).Value = 3181.258  # H122: was 3137.1428
$ws.Cells.Item(122, 9).Value = 1662.76  # I122: was 1625
$ws.Cells.Item(122, 10).Value = 9508.333000000001  # J122: was 7505.5557
$ws.Cells.Item(122, 11).Value = 4988.28  # K122: was 4875
$ws.Cells.Item(122, 12).Value = 28524.999  # L122: was 22516.6671
$ws.Cells.Item(122, 13).Value = -2538.28  # M122: was -2425
$ws.Cells.Item(122, 14).Value = -33424.999  # N122: was -27416.6671
$ws.Cells.Item(132, 8).Value = 17743206  # H132: was 20227124
$ws.Cells.Item(132, 9).Value = 21284410  # I132: was 24399054
$ws.Cells.Item(132, 10).Value = 1099555.6  # J132: was 1221664.5
$ws.Cells.Item(132, 11).Value = 63853230  # K132: was 73197162
$ws.Cells.Item(132, 12).Value = 3298666.8  # L132: was 3664993.5
$ws.Cells.Item(132, 13).Value = -63850700  # M132: was -73194632
$ws.Cells.Item(132, 14).Value = -3303726.8  # N132: was -3670053.5
$ws.Cells.Item(136, 8).Value = 5980.9697  # H136: was 5653.657
$ws.Cells.Item(136, 9).Value = 7573.391  # I136: was 6987.76
$ws.Cells.Item(136, 11).Value = 22720.173  # K136: was 20963.28
$ws.Cells.Item(136, 13).Value = -20170.173  # M136: was -18413.28
